$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104, shifting rows 104:145 down to 105:146
$ws.Rows.Item(104).Insert()

# Populate the new row 104 with the new data point
$ws.Cells.Item(104, 1).Value = 9
$ws.Cells.Item(104, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(104, 3).Value = "Metropolitana"
$ws.Cells.Item(104, 4).Value = 44875
$ws.Cells.Item(104, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(104, 5).Value = 13
$ws.Cells.Item(104, 6).Value = 100112022
$ws.Cells.Item(104, 7).Value = "Arveja Verde"
$ws.Cells.Item(104, 8).Value = "Perfection"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 32
$ws.Cells.Item(104, 11).Value = 20000
$ws.Cells.Item(104, 12).Value = 20000
$ws.Cells.Item(104, 13).Value = 20000
$ws.Cells.Item(104, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(104, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(104, 16).Value = 800
$ws.Cells.Item(104, 17).Value = 25
$ws.Cells.Item(104, 18).Value = "Hortaliza"
